# Fruta / hortaliza, semanal
# Inserts two new weekly price records (Zapallo italiano, Terminal
# Hortofrutícola Agro Chillán) into the existing daily table, shifting the
# subsequent rows down, matching the new weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row above the current row 191 -----------------------
$ws.Rows("191:191").Insert()

$ws.Cells.Item(191, 1).Value  = 7
$ws.Cells.Item(191, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(191, 3).Value  = "Ñuble"
$ws.Cells.Item(191, 4).Value  = 44782
$ws.Cells.Item(191, 5).Value  = 16
$ws.Cells.Item(191, 6).Value  = 100112032
$ws.Cells.Item(191, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(191, 8).Value  = "Sin especificar"
$ws.Cells.Item(191, 9).Value  = "Primera"
$ws.Cells.Item(191, 10).Value = 120
$ws.Cells.Item(191, 11).Value = 20000
$ws.Cells.Item(191, 12).Value = 21000
$ws.Cells.Item(191, 13).Value = 20500
$ws.Cells.Item(191, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(191, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(191, 16).Value = 410
$ws.Cells.Item(191, 17).Value = 50
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# --- Insert second new row above the (now shifted) row 220 ----------------
$ws.Rows("220:220").Insert()

$ws.Cells.Item(220, 1).Value  = 7
$ws.Cells.Item(220, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(220, 3).Value  = "Ñuble"
$ws.Cells.Item(220, 4).Value  = 44783
$ws.Cells.Item(220, 5).Value  = 16
$ws.Cells.Item(220, 6).Value  = 100112032
$ws.Cells.Item(220, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(220, 8).Value  = "Sin especificar"
$ws.Cells.Item(220, 9).Value  = "Primera"
$ws.Cells.Item(220, 10).Value = 80
$ws.Cells.Item(220, 11).Value = 20000
$ws.Cells.Item(220, 12).Value = 21000
$ws.Cells.Item(220, 13).Value = 20500
$ws.Cells.Item(220, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 410
$ws.Cells.Item(220, 17).Value = 50
$ws.Cells.Item(220, 18).Value = "Hortaliza"
